# Generate Report for Handback
# - Status for the f346b3dd... row changes from "Ready for handoff" to
#   "Handback transform failed" (Overview + per-locale sheets all share
#   this status string).
# - The per-locale sheets (zh-cn / de-de) get an Error Detail message in
#   column P for that same row explaining the handback file name mismatch.
# - Column P is widened to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Row 3 on every sheet is the f346b3dd-9b9c-4670-9cea-fb52f7e180db.md entry.
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Error Detail (column P) for that row, per locale.
$zhcn.Range("P3").Value = "Handback file name: ur5md2ua.njp is different with handoff file name: f346b3dd-9b9c-4670-9cea-fb52f7e180db.7acf262da6f3d1b913ace542992d2c6dd0608311.zh-cn."
$dede.Range("P3").Value = "Handback file name: ur5md2ua.njp is different with handoff file name: f346b3dd-9b9c-4670-9cea-fb52f7e180db.7acf262da6f3d1b913ace542992d2c6dd0608311.de-de."

# Widen column P (Error Detail) on both locale sheets to fit the new text.
$zhcn.Columns.Item(16).ColumnWidth = 39.15
$dede.Columns.Item(16).ColumnWidth = 39.15
